{"js": "// Update the \"three-digit number multiplied by one-digit number\" practice\n// table: each cell's multiplication expression is replaced with a new one,\n// keeping the existing run formatting (font, size) intact.\nconst replacements = [\n  { from: \"475\u00d74=\", to: \"476\u00d72=\" },\n  { from: \"430\u00d72=\", to: \"972\u00d76=\" },\n  { from: \"629\u00d77=\", to: \"606\u00d78=\" },\n  { from: \"581\u00d75=\", to: \"333\u00d76=\" },\n  { from: \"221\u00d77=\", to: \"211\u00d78=\" },\n  { from: \"302\u00d72=\", to: \"950\u00d73=\" },\n  { from: \"649\u00d77=\", to: \"585\u00d75=\" },\n  { from: \"948\u00d78=\", to: \"911\u00d74=\" },\n  { from: \"361\u00d79=\", to: \"785\u00d79=\" },\n  { from: \"412\u00d79=\", to: \"763\u00d77=\" },\n  { from: \"441\u00d76=\", to: \"292\u00d72=\" },\n  { from: \"747\u00d73=\", to: \"416\u00d78=\" },\n  { from: \"950\u00d78=\", to: \"702\u00d73=\" },\n  { from: \"127\u00d79=\", to: \"649\u00d78=\" },\n  { from: \"973\u00d79=\", to: \"240\u00d77=\" },\n  { from: \"725\u00d78=\", to: \"332\u00d75=\" },\n  { from: \"666\u00d72=\", to: \"602\u00d78=\" },\n  { from: \"814\u00d72=\", to: \"866\u00d72=\" },\n  { from: \"561\u00d74=\", to: \"111\u00d72=\" },\n  { from: \"648\u00d77=\", to: \"165\u00d75=\" },\n  { from: \"934\u00d77=\", to: \"427\u00d78=\" },\n  { from: \"959\u00d73=\", to: \"971\u00d77=\" },\n  { from: \"274\u00d79=\", to: \"909\u00d73=\" },\n  { from: \"478\u00d73=\", to: \"383\u00d74=\" },\n  { from: \"811\u00d74=\", to: \"889\u00d74=\" },\n];\n\nconst body = context.document.body;\n\nfor (const { from, to } of replacements) {\n  const results = body.search(from, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(to, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the \"three-digit number multiplied by one-digit number\" practice\n# table: each cell's multiplication expression is replaced with a new one.\n# Word's Find/Replace keeps the existing run formatting (font, size) intact.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"475\u00d74=\", \"476\u00d72=\"),\n  @(\"430\u00d72=\", \"972\u00d76=\"),\n  @(\"629\u00d77=\", \"606\u00d78=\"),\n  @(\"581\u00d75=\", \"333\u00d76=\"),\n  @(\"221\u00d77=\", \"211\u00d78=\"),\n  @(\"302\u00d72=\", \"950\u00d73=\"),\n  @(\"649\u00d77=\", \"585\u00d75=\"),\n  @(\"948\u00d78=\", \"911\u00d74=\"),\n  @(\"361\u00d79=\", \"785\u00d79=\"),\n  @(\"412\u00d79=\", \"763\u00d77=\"),\n  @(\"441\u00d76=\", \"292\u00d72=\"),\n  @(\"747\u00d73=\", \"416\u00d78=\"),\n  @(\"950\u00d78=\", \"702\u00d73=\"),\n  @(\"127\u00d79=\", \"649\u00d78=\"),\n  @(\"973\u00d79=\", \"240\u00d77=\"),\n  @(\"725\u00d78=\", \"332\u00d75=\"),\n  @(\"666\u00d72=\", \"602\u00d78=\"),\n  @(\"814\u00d72=\", \"866\u00d72=\"),\n  @(\"561\u00d74=\", \"111\u00d72=\"),\n  @(\"648\u00d77=\", \"165\u00d75=\"),\n  @(\"934\u00d77=\", \"427\u00d78=\"),\n  @(\"959\u00d73=\", \"971\u00d77=\"),\n  @(\"274\u00d79=\", \"909\u00d73=\"),\n  @(\"478\u00d73=\", \"383\u00d74=\"),\n  @(\"811\u00d74=\", \"889\u00d74=\")\n)\n\nforeach ($pair in $pairs) {\n  $rng = $d.Content\n  $rng.Find.ClearFormatting()\n  $rng.Find.Execute($pair[0], $true, $true, $false, $false, $false, $true, 1, $false, $pair[1], 2) | Out-Null\n}\n"}
